$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, pushing existing rows 240:292 down to 241:293
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new weekly record
$ws.Range("A240").Value = 2
$ws.Range("B240").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C240").Value = "Coquimbo"
$ws.Range("D240").Value = 45211
$ws.Range("E240").Value = 4
$ws.Range("F240").Value = 100112031
$ws.Range("G240").Value = "Poroto verde"
$ws.Range("H240").Value = "Magnum"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 500
$ws.Range("K240").Value = 27000
$ws.Range("L240").Value = 29000
$ws.Range("M240").Value = 28000
$ws.Range("N240").Value = "`$/caja 25 kilos"
$ws.Range("O240").Value = "Provincia de Limarí"
$ws.Range("P240").Value = 1120
$ws.Range("Q240").Value = 25
$ws.Range("R240").Value = "Hortaliza"
